# Apply the balanco consolidation edits described in the commit:
# "Adicionados balancos concatenados em uma unica planilha."
#
# The upstream diff clears a batch of placeholder 0-valued cells (turning
# them into blank/empty cells) on several rows, and nudges a handful of
# other numeric cells by tiny floating point amounts (re-aggregation
# rounding from concatenating balance sheets).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows where only AB:AF (the last 5 columns) become blank ---
$rowsAB_AF = @(57, 58, 71, 72, 73, 77, 78)
$colsAB_AF = @("AB", "AC", "AD", "AE", "AF")

foreach ($r in $rowsAB_AF) {
    foreach ($c in $colsAB_AF) {
        $ws.Range("$c$r").Value = ""
    }
}

# --- Rows where B:AB (everything except the trailing AC:AF) become blank ---
$rowsB_AB = @(64, 79)
$colsB_AB = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

foreach ($r in $rowsB_AB) {
    foreach ($c in $colsB_AB) {
        $ws.Range("$c$r").Value = ""
    }
}

# --- Small floating point value adjustments ---
$ws.Range("AF68").Value = 155224.992
$ws.Range("AB70").Value = -96891.016
$ws.Range("X74").Value = -52076.008
$ws.Range("AB74").Value = -86303.992
$ws.Range("X80").Value = -49567
$ws.Range("AB80").Value = -80592.984
